$d = $word.ActiveDocument

# The document currently ends with a single paragraph:
#   "My first Git and Github Class."
# We need to add a new paragraph right after it, inheriting the same
# paragraph formatting (the tab stop at 2160 twips), containing:
#   "I already pushed the project to GitHub, now I'm just making changes."
# (with a curly right single quote, U+2019)

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)             # wdCollapseEnd: collapse to end of the last paragraph
$r.InsertParagraphAfter()  # insert a new paragraph after it, copying pPr (tabs, etc.)

$newPara = $d.Paragraphs.Last
$apost = [char]0x2019
$newPara.Range.Text = "I already pushed the project to GitHub, now I" + $apost + "m just making changes."
